# Auto-generated Excel COM-interop script
# Applies numeric value updates to the FFXIV "Leve profits" sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
# as captured in the upstream commit diff ("chore: update Sheets via scheduled runner").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 82
$ws.Range("H82").Value = 1005
$ws.Range("I82").Value = 1000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2594

# Row 85
$ws.Range("H85").Value = 1005
$ws.Range("I85").Value = 1000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1596

# Row 98
$ws.Range("H98").Value = 2372.9023
$ws.Range("I98").Value = 2372.9023
$ws.Range("K98").Value = 2372.9023
$ws.Range("M98").Value = -874.9023000000002

# Row 122
$ws.Range("H122").Value = 2372.9023
$ws.Range("I122").Value = 2372.9023
$ws.Range("K122").Value = 7118.706900000001
$ws.Range("M122").Value = -4668.706900000001

# Row 132
$ws.Range("H132").Value = 1695.8334
$ws.Range("I132").Value = 1695.8334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5087.5002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3728
$ws.Range("I45").Value = 1856.5625
$ws.Range("J45").Value = 6722.3
$ws.Range("K45").Value = 1856.5625
$ws.Range("L45").Value = 6722.3
$ws.Range("M45").Value = -1479.5625
$ws.Range("N45").Value = -7476.3

# Row 61
$ws.Range("H61").Value = 32263556
$ws.Range("I61").Value = 2123.8572
$ws.Range("J61").Value = 100012570
$ws.Range("K61").Value = 2123.8572
$ws.Range("L61").Value = 100012570
$ws.Range("M61").Value = -1911.8572
$ws.Range("N61").Value = -100012994

# Row 122
$ws.Range("H122").Value = 2822.24
$ws.Range("I122").Value = 1939
$ws.Range("K122").Value = 5817
$ws.Range("M122").Value = -3367

# Row 136
$ws.Range("H136").Value = 32263556
$ws.Range("I136").Value = 2123.8572
$ws.Range("J136").Value = 100012570
$ws.Range("K136").Value = 6371.571599999999
$ws.Range("L136").Value = 300037710
$ws.Range("M136").Value = -3821.571599999999
$ws.Range("N136").Value = -300042810

$ws = $wb.Worksheets.Item("BSM")
# Row 60
$ws.Range("H60").Value = 178996
$ws.Range("J60").Value = 178996
$ws.Range("L60").Value = 178996
$ws.Range("N60").Value = -180194

# Row 86
$ws.Range("H86").Value = 8369173.5
$ws.Range("I86").Value = 11410980
$ws.Range("K86").Value = 11410980
$ws.Range("M86").Value = -11409857

# Row 89
$ws.Range("H89").Value = 8369173.5
$ws.Range("I89").Value = 11410980
$ws.Range("K89").Value = 57054900
$ws.Range("M89").Value = -57049284

# Row 134
$ws.Range("H134").Value = 7358498.5
$ws.Range("I134").Value = 13158861
$ws.Range("K134").Value = 39476583
$ws.Range("M134").Value = -39474048

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 311.33334
$ws.Range("I7").Value = 361.07144
$ws.Range("K7").Value = 361.07144
$ws.Range("M7").Value = -248.07144

# Row 16
$ws.Range("H16").Value = 4678.4585
$ws.Range("I16").Value = 1999.1818
$ws.Range("K16").Value = 1999.1818
$ws.Range("M16").Value = -1712.1818

# Row 31
$ws.Range("H31").Value = 8559.269
$ws.Range("I31").Value = 2341.5557
$ws.Range("J31").Value = 10308
$ws.Range("K31").Value = 2341.5557
$ws.Range("L31").Value = 10308
$ws.Range("M31").Value = -2046.5557
$ws.Range("N31").Value = -10898

# Row 34
$ws.Range("H34").Value = 8559.269
$ws.Range("I34").Value = 2341.5557
$ws.Range("J34").Value = 10308
$ws.Range("K34").Value = 2341.5557
$ws.Range("L34").Value = 10308
$ws.Range("M34").Value = -2139.5557
$ws.Range("N34").Value = -10712

# Row 58
$ws.Range("H58").Value = 6316.298
$ws.Range("I58").Value = 3459.4211
$ws.Range("K58").Value = 3459.4211
$ws.Range("M58").Value = -3256.4211

# Row 62
$ws.Range("H62").Value = 5214060
$ws.Range("J62").Value = 6461
$ws.Range("L62").Value = 6461
$ws.Range("N62").Value = -7709

# Row 65
$ws.Range("H65").Value = 5214060
$ws.Range("J65").Value = 6461
$ws.Range("L65").Value = 32305
$ws.Range("N65").Value = -38545

# Row 99
$ws.Range("H99").Value = 3979.96
$ws.Range("I99").Value = 2974.9375
$ws.Range("K99").Value = 2974.9375
$ws.Range("M99").Value = -1476.9375

# Row 110
$ws.Range("H110").Value = 69500
$ws.Range("J110").Value = 69500
$ws.Range("L110").Value = 69500
$ws.Range("N110").Value = -77680

# Row 113
$ws.Range("H113").Value = 4678.4585
$ws.Range("I113").Value = 1999.1818
$ws.Range("K113").Value = 1999.1818
$ws.Range("M113").Value = 170.8181999999999

# Row 126
$ws.Range("H126").Value = 3979.96
$ws.Range("I126").Value = 2974.9375
$ws.Range("K126").Value = 8924.8125
$ws.Range("M126").Value = -6454.8125

# Row 133
$ws.Range("H133").Value = 51000
$ws.Range("J133").Value = 51000
$ws.Range("L133").Value = 51000
$ws.Range("N133").Value = -56060

# Row 136
$ws.Range("H136").Value = 6316.298
$ws.Range("I136").Value = 3459.4211
$ws.Range("K136").Value = 10378.2633
$ws.Range("M136").Value = -7828.263300000001

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 1899.75
$ws.Range("I3").Value = 1899.75
$ws.Range("K3").Value = 5699.25
$ws.Range("M3").Value = -5587.25

# Row 12
$ws.Range("H12").Value = 917.4138
$ws.Range("J12").Value = 464.5
$ws.Range("L12").Value = 1393.5
$ws.Range("N12").Value = -1739.5

# Row 34
$ws.Range("H34").Value = 5736.9473
$ws.Range("I34").Value = 2000
$ws.Range("J34").Value = 5944.5557
$ws.Range("K34").Value = 6000
$ws.Range("L34").Value = 17833.6671
$ws.Range("M34").Value = -5916
$ws.Range("N34").Value = -18001.6671

# Row 113
$ws.Range("H113").Value = 6782.353
$ws.Range("I113").Value = 1854.4
$ws.Range("J113").Value = 8835.666999999999
$ws.Range("K113").Value = 5563.200000000001
$ws.Range("L113").Value = 26507.001
$ws.Range("M113").Value = -3393.200000000001
$ws.Range("N113").Value = -30847.001

# Row 131
$ws.Range("H131").Value = 31950.455
$ws.Range("J131").Value = 32902.03
$ws.Range("L131").Value = 98706.09
$ws.Range("N131").Value = -108786.09

# Row 132
$ws.Range("H132").Value = 10670.143
$ws.Range("I132").Value = 3539
$ws.Range("J132").Value = 16018.5
$ws.Range("K132").Value = 31851
$ws.Range("L132").Value = 144166.5
$ws.Range("M132").Value = -29321
$ws.Range("N132").Value = -149226.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6863.343
$ws.Range("I70").Value = 5975.7085
$ws.Range("K70").Value = 5975.7085
$ws.Range("M70").Value = -5705.7085

# Row 73
$ws.Range("H73").Value = 6863.343
$ws.Range("I73").Value = 5975.7085
$ws.Range("K73").Value = 5975.7085
$ws.Range("M73").Value = -5039.7085

# Row 80
$ws.Range("H80").Value = 2871.3572
$ws.Range("I80").Value = 2750
$ws.Range("J80").Value = 3174.75
$ws.Range("K80").Value = 2750
$ws.Range("L80").Value = 3174.75
$ws.Range("M80").Value = -1752
$ws.Range("N80").Value = -5170.75

# Row 83
$ws.Range("H83").Value = 2871.3572
$ws.Range("I83").Value = 2750
$ws.Range("J83").Value = 3174.75
$ws.Range("K83").Value = 13750
$ws.Range("L83").Value = 15873.75
$ws.Range("M83").Value = -8758
$ws.Range("N83").Value = -25857.75

# Row 101
$ws.Range("H101").Value = 38492
$ws.Range("J101").Value = 38492
$ws.Range("L101").Value = 38492
$ws.Range("N101").Value = -44982

# Row 132
$ws.Range("H132").Value = 5510.737
$ws.Range("I132").Value = 2091.1
$ws.Range("J132").Value = 9310.333000000001
$ws.Range("K132").Value = 6273.299999999999
$ws.Range("L132").Value = 27930.999
$ws.Range("M132").Value = -3743.299999999999
$ws.Range("N132").Value = -32990.999

# Row 136
$ws.Range("H136").Value = 33968.406
$ws.Range("J136").Value = 34658.566
$ws.Range("L136").Value = 103975.698
$ws.Range("N136").Value = -109075.698

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1835.75
$ws.Range("I16").Value = 1812.2858
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1812.2858
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1642.2858
$ws.Range("N16").Value = -2340

# Row 22
$ws.Range("H22").Value = 814.71875
$ws.Range("I22").Value = 303.15384
$ws.Range("K22").Value = 303.15384
$ws.Range("M22").Value = -8.153840000000002

# Row 27
$ws.Range("H27").Value = 814.71875
$ws.Range("I27").Value = 303.15384
$ws.Range("K27").Value = 303.15384
$ws.Range("M27").Value = -196.15384

# Row 46
$ws.Range("H46").Value = 3390.6191
$ws.Range("I46").Value = 2247.2307
$ws.Range("K46").Value = 2247.2307
$ws.Range("M46").Value = -2059.2307

# Row 55
$ws.Range("H55").Value = 321.5357
$ws.Range("I55").Value = 154.38461
$ws.Range("J55").Value = 466.4
$ws.Range("K55").Value = 154.38461
$ws.Range("L55").Value = 466.4
$ws.Range("M55").Value = 18.61538999999999
$ws.Range("N55").Value = -812.4

# Row 132
$ws.Range("H132").Value = 7817564
$ws.Range("I132").Value = 12197663
$ws.Range("K132").Value = 36592989
$ws.Range("M132").Value = -36590459

# Row 136
$ws.Range("H136").Value = 10282.086
$ws.Range("I136").Value = 3344.2415
$ws.Range("J136").Value = 17219.932
$ws.Range("K136").Value = 10032.7245
$ws.Range("L136").Value = 51659.796
$ws.Range("M136").Value = -7482.7245
$ws.Range("N136").Value = -56759.796

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 7282.385
$ws.Range("I132").Value = 7244.0586
$ws.Range("J132").Value = 7354.778
$ws.Range("K132").Value = 21732.1758
$ws.Range("L132").Value = 22064.334
$ws.Range("M132").Value = -19202.1758
$ws.Range("N132").Value = -27124.334

# Row 136
$ws.Range("H136").Value = 26583950
$ws.Range("I136").Value = 55558416
$ws.Range("K136").Value = 166675248
$ws.Range("M136").Value = -166672698
